$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.200.73"
$ws.Range("E2").Value = "  -0.70%  "

$ws.Range("D3").Value = "1.859.84"
$ws.Range("E3").Value = "  -0.99%  "

$ws.Range("D4").Value = "'0.9997"
$ws.Range("E4").Value = "  -0.13%  "

$ws.Range("D5").Value = "'241.83"
$ws.Range("E5").Value = "  -0.77%  "

$ws.Range("D6").Value = "'0.7010"
$ws.Range("E6").Value = "  -2.17%  "

$ws.Range("D7").Value = "'0.9998"
$ws.Range("E7").Value = "  -0.14%  "

$ws.Range("D8").Value = "'0.07812"
$ws.Range("E8").Value = "  -2.16%  "

$ws.Range("E9").Value = "  -1.34%  "

$ws.Range("D10").Value = "'23.91"
$ws.Range("E10").Value = "  -4.17%  "

$ws.Range("D11").Value = "'0.07804"
$ws.Range("E11").Value = "  -3.52%  "

$ws.Range("D12").Value = "1.862.81"
$ws.Range("E12").Value = "  -1.13%  "

$ws.Range("D13").Value = "'92.66"
$ws.Range("E13").Value = "  -2.19%  "

$ws.Range("E14").Value = "  -2.00%  "

$ws.Range("D15").Value = "'0.6909"
$ws.Range("E15").Value = "  -2.41%  "

$ws.Range("D16").Value = "'6.563"
$ws.Range("E16").Value = "  +2.49%  "

$ws.Range("D17").Value = "'0.000008435"
$ws.Range("E17").Value = "  -0.13%  "

$ws.Range("D18").Value = "29.228.53"
$ws.Range("E18").Value = "  -0.62%  "

$ws.Range("D19").Value = "'249.75"
$ws.Range("E19").Value = "  -1.32%  "

$ws.Range("D20").Value = "2.112.09"
$ws.Range("E20").Value = "  -1.08%  "

$ws.Range("D21").Value = "'12.93"
$ws.Range("E21").Value = "  -3.19%  "

$ws.Range("D22").Value = "'0.9993"
$ws.Range("E22").Value = "  -0.14%  "

$ws.Range("D23").Value = "'7.595"
$ws.Range("E23").Value = "  -1.12%  "

$ws.Range("D24").Value = "'0.9999"
$ws.Range("E24").Value = "  -0.11%  "

$ws.Range("D25").Value = "'0.1532"
$ws.Range("E25").Value = "  -3.12%  "

$ws.Range("D26").Value = "'160.88"
$ws.Range("E26").Value = "  -0.89%  "

$ws.Range("D27").Value = "'8.896"
$ws.Range("E27").Value = "  -1.97%  "

$ws.Range("E28").Value = "  -2.19%  "

$ws.Range("D29").Value = "'1.572"
$ws.Range("E29").Value = "  +4.28%  "

$ws.Range("D30").Value = "'4.278"
$ws.Range("E30").Value = "  -3.24%  "

$ws.Range("D31").Value = "'4.250"
$ws.Range("E31").Value = "  -1.61%  "

$ws.Range("D32").Value = "'1.214"
$ws.Range("E32").Value = "  -0.68%  "

$ws.Range("D33").Value = "'0.05225"
$ws.Range("E33").Value = "  -1.75%  "

$ws.Range("D34").Value = "'0.7584"
$ws.Range("E34").Value = "  -0.07%  "

$ws.Range("D35").Value = "'1.875"
$ws.Range("E35").Value = "  -3.71%  "

$ws.Range("D37").Value = "'2.710"
$ws.Range("E37").Value = "  +0.30%  "

$ws.Range("E38").Value = "  -1.34%  "

$ws.Range("D39").Value = "1.220.43"
$ws.Range("E39").Value = "  -4.52%  "

$ws.Range("D40").Value = "'2.721"
$ws.Range("E40").Value = "  -1.34%  "

$ws.Range("D41").Value = "'0.9002"
$ws.Range("E41").Value = "  -0.82%  "

$ws.Range("E42").Value = "  -1.18%  "

$ws.Range("D43").Value = "'5.815"
$ws.Range("E43").Value = "  -9.30%  "

$ws.Range("D44").Value = "'0.9992"
$ws.Range("E44").Value = "  -0.16%  "

$ws.Range("D45").Value = "'66.51"
$ws.Range("E45").Value = "  -10.41%  "

$ws.Range("D46").Value = "2.010.47"
$ws.Range("E46").Value = "  -0.99%  "

$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "'0.00000000124"
$ws.Range("E47").Value = "  -4.56%  "

$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").Value = "'0.5178"
$ws.Range("E48").Value = "  -0.58%  "

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'9.528"
$ws.Range("E49").Value = "  -0.02%  "

$ws.Range("D50").Value = "'1.769"
$ws.Range("E50").Value = "  -2.10%  "

$ws.Range("D51").Value = "'7.047"
$ws.Range("E51").Value = "  -0.81%  "
